$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from H1 onto the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row (I1:J1) - new columns "I0" and "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-30: I = 1, J = same value as H in that row
for ($r = 2; $r -le 30; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}

# Row 31 special values
$ws.Cells.Item(31, 9).Value = 4
$ws.Cells.Item(31, 10).Value = 6

# Row 32 special values
$ws.Cells.Item(32, 9).Value = 3
$ws.Cells.Item(32, 10).Value = 4
